$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells (which look numeric) are written as text, not converted to numbers.
function Set-TextCell($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell $ws 'D2' '29.221.03'
Set-TextCell $ws 'E2' '  -2.44%  '

Set-TextCell $ws 'D3' '1.852.71'
Set-TextCell $ws 'E3' '  -1.35%  '

Set-TextCell $ws 'D4' '0.9999'
Set-TextCell $ws 'E4' '  -0.14%  '

Set-TextCell $ws 'D5' '0.6972'
Set-TextCell $ws 'E5' '  -5.82%  '

Set-TextCell $ws 'D6' '238.71'
Set-TextCell $ws 'E6' '  -1.80%  '

Set-TextCell $ws 'D7' '1.000'
Set-TextCell $ws 'E7' '  -0.26%  '

Set-TextCell $ws 'D8' '0.3073'
Set-TextCell $ws 'E8' '  -2.52%  '

Set-TextCell $ws 'D9' '0.07608'
Set-TextCell $ws 'E9' '  +5.62%  '

Set-TextCell $ws 'E10' '  -4.08%  '

Set-TextCell $ws 'D11' '0.08090'
Set-TextCell $ws 'E11' '  -3.03%  '

Set-TextCell $ws 'D12' '1.850.57'
Set-TextCell $ws 'E12' '  -1.18%  '

Set-TextCell $ws 'D13' '0.7257'
Set-TextCell $ws 'E13' '  -3.30%  '

Set-TextCell $ws 'E14' '  -3.98%  '

Set-TextCell $ws 'D15' '89.14'
Set-TextCell $ws 'E15' '  -3.60%  '

Set-TextCell $ws 'D16' '29.216.43'
Set-TextCell $ws 'E16' '  -2.56%  '

Set-TextCell $ws 'D17' '5.870'
Set-TextCell $ws 'E17' '  -3.80%  '

Set-TextCell $ws 'D18' '241.70'
Set-TextCell $ws 'E18' '  -2.69%  '

Set-TextCell $ws 'D19' '0.000007721'
Set-TextCell $ws 'E19' '  -1.66%  '

Set-TextCell $ws 'E20' '  -3.25%  '

Set-TextCell $ws 'D21' '0.9994'
Set-TextCell $ws 'E21' '  -0.11%  '

Set-TextCell $ws 'D22' '2.095.42'
Set-TextCell $ws 'E22' '  -2.61%  '

Set-TextCell $ws 'D23' '0.9998'
Set-TextCell $ws 'E23' '  -0.10%  '

Set-TextCell $ws 'D24' '7.630'
Set-TextCell $ws 'E24' '  -4.84%  '

Set-TextCell $ws 'D25' '9.051'
Set-TextCell $ws 'E25' '  -2.45%  '

Set-TextCell $ws 'D26' '161.63'
Set-TextCell $ws 'E26' '  -1.98%  '

Set-TextCell $ws 'D27' '0.1462'
Set-TextCell $ws 'E27' '  -5.49%  '

Set-TextCell $ws 'D28' '18.07'
Set-TextCell $ws 'E28' '  -3.35%  '

Set-TextCell $ws 'D29' '1.936'
Set-TextCell $ws 'E29' '  -4.86%  '

Set-TextCell $ws 'D30' '1.393'
Set-TextCell $ws 'E30' '  -7.78%  '

Set-TextCell $ws 'E31' '  -3.38%  '

Set-TextCell $ws 'D32' '1.503'
Set-TextCell $ws 'E32' '  -2.22%  '

Set-TextCell $ws 'D33' '4.049'
Set-TextCell $ws 'E33' '  -5.14%  '

Set-TextCell $ws 'D34' '0.05262'
Set-TextCell $ws 'E34' '  -1.10%  '

Set-TextCell $ws 'D35' '1.194'
Set-TextCell $ws 'E35' '  -3.32%  '

Set-TextCell $ws 'D36' '0.7107'
Set-TextCell $ws 'E36' '  -5.18%  '

Set-TextCell $ws 'D37' '1.003'
Set-TextCell $ws 'E37' '  +0.16%  '

Set-TextCell $ws 'D38' '2.665'
Set-TextCell $ws 'E38' '  -1.10%  '

Set-TextCell $ws 'E39' '  -5.34%  '

Set-TextCell $ws 'D40' '2.689'
Set-TextCell $ws 'E40' '  -2.49%  '

Set-TextCell $ws 'D41' '0.9219'
Set-TextCell $ws 'E41' '  +7.35%  '

Set-TextCell $ws 'D42' '5.965'
Set-TextCell $ws 'E42' '  -2.89%  '

Set-TextCell $ws 'D43' '0.4303'
Set-TextCell $ws 'E43' '  -5.36%  '

Set-TextCell $ws 'B44' 'Maker'
Set-TextCell $ws 'C44' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell $ws 'D44' '1.047.82'
Set-TextCell $ws 'E44' '  -5.36%  '

Set-TextCell $ws 'B45' 'Aave'
Set-TextCell $ws 'C45' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell $ws 'D45' '69.63'
Set-TextCell $ws 'E45' '  -3.59%  '

Set-TextCell $ws 'E46' '  -0.24%  '

Set-TextCell $ws 'D47' '102.31'
Set-TextCell $ws 'E47' '  -1.92%  '

Set-TextCell $ws 'D48' '7.247'
Set-TextCell $ws 'E48' '  -4.85%  '

Set-TextCell $ws 'D49' '1.742'
Set-TextCell $ws 'E49' '  -6.18%  '

Set-TextCell $ws 'D50' '9.270'
Set-TextCell $ws 'E50' '  -2.43%  '

Set-TextCell $ws 'D51' '1.998.27'
Set-TextCell $ws 'E51' '  -2.05%  '
